# Fix data that was incorrectly pulled in from "extra files" (other tickers)
# mixed into the SALESFORCECOM_INC.xlsx dataset.
#
# For every data row (2-44) on the active sheet:
#   - columns D,E,F,G (open/close/high/low price) and H (shares_outstanding)
#     are corrected to the values that actually belong to CRM
#   - column I (fixed_ticker) is corrected to "CRM" for every row
#   - the now-unused ticker strings that leaked in from the other files
#     (QCOM, MRVL, NOW, ... PTC) are removed from the workbook

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, open_price(D), close_price(E), high_price(F), low_price(G), shares_outstanding(H)
$data = @(
  "2,56.2879002583733,68.65777587890625,70.2609149665568,54.87278575352673,952000000",
  "3,72.60623971603434,71.99269104003906,75.38698591615947,68.56870966426919,952000000",
  "4,72.80415957253398,68.63798522949219,73.36822575318394,63.49211858321099,952000000",
  "5,76.90107357485175,78.86046600341797,82.0370504671741,73.1703208706824,952000000",
  "6,66.88640070712104,67.04473876953125,70.10256898757827,52.05244512140504,952000000",
  "7,75.50575413074584,82.83861541748047,83.60060428970398,72.40833119946025,952000000",
  "8,80.92870334345369,78.59326934814453,80.94849937984105,75.08022788696286,952000000",
  "9,75.19896683055654,71.25050354003906,79.53337729455961,71.25050354003906,952000000",
  "10,78.08858589770188,80.50318145751953,82.26465124746039,76.64377719923537,952000000",
  "11,85.35217530942711,88.70688629150391,90.85429573077296,85.21363325205336,952000000",
  "12,90.41887194645859,94.49597930908205,94.6642154046076,86.35165505123808,952000000",
  "13,102.5908269748701,103.2340621948242,108.0533697858935,100.2455027413109,952000000",
  "14,112.1205903160145,115.0398788452148,116.9201020309596,101.3043674266762,952000000",
  "15,119.6216903362453,127.9837341308594,131.170216573424,119.166479637304,952000000",
  "16,136.563472887765,151.0906524658203,153.3568228994231,135.3264849087091,952000000",
  "17,136.9593108587654,141.27392578125,142.1348797320345,112.4174724762589,952000000",
  "18,150.8135511125558,161.9464416503906,162.8865644868783,149.4974093413715,952000000",
  "19,164.2819209818063,149.8339080810547,164.4897416402893,149.6260874225716,952000000",
  "20,153.11931098708,154.4453735351562,156.770901803147,136.4348170429022,952000000",
  "21,156.2662404791025,161.1943969726562,164.6084811995161,153.8516372276565,952000000",
  "22,181.3523227347343,168.6261901855469,193.6826251825274,161.4219736205274,952000000",
  "23,156.8401789823443,172.9704895019531,182.8762985046751,153.4656722869126,952000000",
  "24,195.7409707736601,269.8117980957031,275.3831964992401,189.7242597405219,952000000",
  "25,232.6724700721815,243.2413024902344,268.0998214754197,226.6062917081189,952000000",
  "26,224.8448264972343,214.2463073730469,248.6147754548341,214.0879692955729,952000000",
  "27,228.367761142473,235.6214599609375,240.9652416201309,206.7353162599697,952000000",
  "28,239.7875934264296,262.5085754394531,272.3550097827603,236.3933063665365,952000000",
  "29,298.3615171679053,281.9936828613281,308.5048191611757,280.0936861316079,952000000",
  "30,229.9511048225611,208.3384552001953,232.049046632756,182.5200465898714,952000000",
  "31,173.2673678392397,158.5719604492188,184.0440059652352,152.9411890365635,952000000",
  "32,179.412700264181,154.4948120117188,192.3466358925259,154.346379526675,952000000",
  "33,163.2824178567729,158.5818634033203,163.9751281246598,134.6238729029025,952000000",
  "34,165.9839986090795,161.9068908691406,176.9783488912497,157.9980116903209,952000000",
  "35,195.7013735792695,221.0546875,222.6578189645609,188.7742414337054,952000000",
  "36,221.5593909943528,219.1546936035156,226.4083778474535,198.1159923258932,952000000",
  "37,198.8285021954522,249.2778167724609,249.8715770549482,198.4029855696322,952000000",
  "38,279.0743998164025,305.6053161621094,307.3074128669396,273.54259307568,952000000",
  "39,266.5673505552567,232.3026275634766,287.3758948652332,210.067209224525,952000000",
  "40,257.9313684090956,250.984390258789,269.8901132539759,230.4511204176801,952000000",
  "41,288.4838176559762,328.0045776367188,346.761039837357,288.2551905501395,952000000",
  "42,334.3296604260772,296.3957824707031,346.340737115616,290.7236000178734,952000000",
  "43,272.4926372448254,264.4892883300781,295.0674606106557,253.6553684235463,952000000",
  "44,254.5467611476779,255.8045654296875,257.1921638417244,226.0863097974658,952000000"
)

foreach ($line in $data) {
  $parts = $line.Split(",")
  $row = [int]$parts[0]
  $openPrice = [double]$parts[1]
  $closePrice = [double]$parts[2]
  $highPrice = [double]$parts[3]
  $lowPrice = [double]$parts[4]
  $shares = [double]$parts[5]

  $ws.Cells.Item($row, 4).Value = $openPrice
  $ws.Cells.Item($row, 5).Value = $closePrice
  $ws.Cells.Item($row, 6).Value = $highPrice
  $ws.Cells.Item($row, 7).Value = $lowPrice
  $ws.Cells.Item($row, 8).Value = $shares
  # column I = fixed_ticker -> always "CRM" after the fix
  $ws.Cells.Item($row, 9).Value = "CRM"
}
